# Auto-generated edit script applying value updates to Ultima Profits sheets
$wb = $excel.ActiveWorkbook

# ALC row 18
$ws = $wb.Worksheets.Item(1)
$ws.Range("H18").Value = 2812.125
$ws.Range("I18").Value = 499.5
$ws.Range("J18").Value = 3583
$ws.Range("K18").Value = 499.5
$ws.Range("L18").Value = 3583
$ws.Range("M18").Value = -215.5
$ws.Range("N18").Value = -4151

# ALC row 106
$ws = $wb.Worksheets.Item(1)
$ws.Range("H106").Value = 101882
$ws.Range("I106").Value = 101882
$ws.Range("K106").Value = 101882
$ws.Range("M106").Value = -101251

# ALC row 125
$ws = $wb.Worksheets.Item(1)
$ws.Range("H125").Value = 5398
$ws.Range("I125").Value = 10000
$ws.Range("J125").Value = 796
$ws.Range("K125").Value = 90000
$ws.Range("L125").Value = 7164
$ws.Range("M125").Value = -87540
$ws.Range("N125").Value = -12084

# ALC row 132
$ws = $wb.Worksheets.Item(1)
$ws.Range("H132").Value = 7579360
$ws.Range("I132").Value = 3280.238
$ws.Range("J132").Value = 20837500
$ws.Range("K132").Value = 9840.714
$ws.Range("L132").Value = 62512500
$ws.Range("M132").Value = -7310.714
$ws.Range("N132").Value = -62517560

# ALC row 135
$ws = $wb.Worksheets.Item(1)
$ws.Range("H135").Value = 12500951
$ws.Range("I135").Value = 916.45715
$ws.Range("J135").Value = 100001200
$ws.Range("K135").Value = 8248.11435
$ws.Range("L135").Value = 900010800
$ws.Range("M135").Value = -5713.11435
$ws.Range("N135").Value = -900015870

# ALC row 137
$ws = $wb.Worksheets.Item(1)
$ws.Range("H137").Value = 13336396
$ws.Range("I137").Value = 2034
$ws.Range("J137").Value = 25003962
$ws.Range("K137").Value = 6102
$ws.Range("L137").Value = 75011886
$ws.Range("M137").Value = -3552
$ws.Range("N137").Value = -75016986

# ALC row 138
$ws = $wb.Worksheets.Item(1)
$ws.Range("H138").Value = 8773936
$ws.Range("I138").Value = 13890069
$ws.Range("J138").Value = 3422.6428
$ws.Range("K138").Value = 41670207
$ws.Range("L138").Value = 10267.9284
$ws.Range("M138").Value = -41665067
$ws.Range("N138").Value = -20547.9284

# ALC row 141
$ws = $wb.Worksheets.Item(1)
$ws.Range("H141").Value = 3707.7097
$ws.Range("I141").Value = 1703
$ws.Range("J141").Value = 4662.3335
$ws.Range("K141").Value = 5109
$ws.Range("L141").Value = 13987.0005
$ws.Range("M141").Value = 71
$ws.Range("N141").Value = -24347.0005

# ARM row 32
$ws = $wb.Worksheets.Item(2)
$ws.Range("H32").Value = 10033.3545
$ws.Range("I32").Value = 9413.671
$ws.Range("J32").Value = 12803.706
$ws.Range("K32").Value = 9413.671
$ws.Range("L32").Value = 12803.706
$ws.Range("M32").Value = -9126.671
$ws.Range("N32").Value = -13377.706

# ARM row 61
$ws = $wb.Worksheets.Item(2)
$ws.Range("H61").Value = 26318250
$ws.Range("I61").Value = 29414338
$ws.Range("J61").Value = 1500
$ws.Range("K61").Value = 29414338
$ws.Range("L61").Value = 1500
$ws.Range("M61").Value = -29414126
$ws.Range("N61").Value = -1924

# ARM row 102
$ws = $wb.Worksheets.Item(2)
$ws.Range("H102").Value = 5000
$ws.Range("I102").Value = 5000
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 5000
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = -3378
$ws.Range("N102").Value = ""

# ARM row 132
$ws = $wb.Worksheets.Item(2)
$ws.Range("H132").Value = 8066024
$ws.Range("I132").Value = 10417789
$ws.Range("J132").Value = 2830.2856
$ws.Range("K132").Value = 31253367
$ws.Range("L132").Value = 8490.856800000001
$ws.Range("M132").Value = -31250837
$ws.Range("N132").Value = -13550.8568

# ARM row 136
$ws = $wb.Worksheets.Item(2)
$ws.Range("H136").Value = 26318250
$ws.Range("I136").Value = 29414338
$ws.Range("J136").Value = 1500
$ws.Range("K136").Value = 88243014
$ws.Range("L136").Value = 4500
$ws.Range("M136").Value = -88240464
$ws.Range("N136").Value = -9600

# BSM row 62
$ws = $wb.Worksheets.Item(3)
$ws.Range("H62").Value = 47681
$ws.Range("J62").Value = 47681
$ws.Range("L62").Value = 47681
$ws.Range("N62").Value = -49053

# BSM row 65
$ws = $wb.Worksheets.Item(3)
$ws.Range("H65").Value = 47681
$ws.Range("J65").Value = 47681
$ws.Range("L65").Value = 143043
$ws.Range("N65").Value = -149907

# BSM row 105
$ws = $wb.Worksheets.Item(3)
$ws.Range("H105").Value = 2668.3386
$ws.Range("I105").Value = 1383.0278
$ws.Range("J105").Value = 4448
$ws.Range("K105").Value = 1383.0278
$ws.Range("L105").Value = 4448
$ws.Range("M105").Value = 363.9721999999999
$ws.Range("N105").Value = -7942

# BSM row 134
$ws = $wb.Worksheets.Item(3)
$ws.Range("H134").Value = 2905.1428
$ws.Range("I134").Value = 2278.72
$ws.Range("J134").Value = 4471.2
$ws.Range("K134").Value = 6836.16
$ws.Range("L134").Value = 13413.6
$ws.Range("M134").Value = -4301.16
$ws.Range("N134").Value = -18483.6

# CRP row 105
$ws = $wb.Worksheets.Item(4)
$ws.Range("H105").Value = 3300
$ws.Range("I105").Value = 3300
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 3300
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = -1553
$ws.Range("N105").Value = ""

# CRP row 134
$ws = $wb.Worksheets.Item(4)
$ws.Range("H134").Value = 682305.1
$ws.Range("I134").Value = 2562.6667
$ws.Range("J134").Value = 1402032.4
$ws.Range("K134").Value = 7688.000100000001
$ws.Range("L134").Value = 4206097.199999999
$ws.Range("M134").Value = -5153.000100000001
$ws.Range("N134").Value = -4211167.199999999

# CUL row 107
$ws = $wb.Worksheets.Item(5)
$ws.Range("H107").Value = 281.55554
$ws.Range("I107").Value = 202.41667
$ws.Range("J107").Value = 439.83334
$ws.Range("K107").Value = 607.25001
$ws.Range("L107").Value = 1319.50002
$ws.Range("M107").Value = 1312.74999
$ws.Range("N107").Value = -5159.500019999999

# CUL row 118
$ws = $wb.Worksheets.Item(5)
$ws.Range("H118").Value = 1617.8
$ws.Range("J118").Value = 1407.9
$ws.Range("L118").Value = 4223.700000000001
$ws.Range("N118").Value = -6709.700000000001

# CUL row 132
$ws = $wb.Worksheets.Item(5)
$ws.Range("H132").Value = 1198.7273
$ws.Range("I132").Value = 757.2632
$ws.Range("J132").Value = 1797.8572
$ws.Range("K132").Value = 6815.3688
$ws.Range("L132").Value = 16180.7148
$ws.Range("M132").Value = -4285.3688
$ws.Range("N132").Value = -21240.7148

# GSM row 2
$ws = $wb.Worksheets.Item(6)
$ws.Range("H2").Value = 63.266666
$ws.Range("I2").Value = 30.75
$ws.Range("J2").Value = 193.33333
$ws.Range("K2").Value = 30.75
$ws.Range("L2").Value = 193.33333
$ws.Range("M2").Value = 82.25
$ws.Range("N2").Value = -419.33333

# GSM row 113
$ws = $wb.Worksheets.Item(6)
$ws.Range("H113").Value = 500855.5
$ws.Range("I113").Value = 1000011
$ws.Range("J113").Value = 1700
$ws.Range("K113").Value = 1000011
$ws.Range("L113").Value = 1700
$ws.Range("M113").Value = -997841
$ws.Range("N113").Value = -6040

# GSM row 122
$ws = $wb.Worksheets.Item(6)
$ws.Range("H122").Value = 2900219.2
$ws.Range("I122").Value = 3704941.5
$ws.Range("J122").Value = 3219.2
$ws.Range("K122").Value = 11114824.5
$ws.Range("L122").Value = 9657.599999999999
$ws.Range("M122").Value = -11112374.5
$ws.Range("N122").Value = -14557.6

# GSM row 132
$ws = $wb.Worksheets.Item(6)
$ws.Range("H132").Value = 2246.3696
$ws.Range("I132").Value = 1988.9445
$ws.Range("K132").Value = 5966.833500000001
$ws.Range("M132").Value = -3436.833500000001

# LTW row 32
$ws = $wb.Worksheets.Item(7)
$ws.Range("H32").Value = 8840
$ws.Range("I32").Value = 3600
$ws.Range("J32").Value = 29800
$ws.Range("K32").Value = 3600
$ws.Range("L32").Value = 29800
$ws.Range("M32").Value = -3283
$ws.Range("N32").Value = -30434

# LTW row 40
$ws = $wb.Worksheets.Item(7)
$ws.Range("H40").Value = 9777.777
$ws.Range("J40").Value = 5000
$ws.Range("L40").Value = 5000
$ws.Range("N40").Value = -5272

# LTW row 132
$ws = $wb.Worksheets.Item(7)
$ws.Range("H132").Value = 8477469
$ws.Range("I132").Value = 2210.8914
$ws.Range("J132").Value = 38466844
$ws.Range("K132").Value = 6632.674199999999
$ws.Range("L132").Value = 115400532
$ws.Range("M132").Value = -4102.674199999999
$ws.Range("N132").Value = -115405592

# LTW row 136
$ws = $wb.Worksheets.Item(7)
$ws.Range("H136").Value = 26317264
$ws.Range("I136").Value = 29413164
$ws.Range("J136").Value = 2111
$ws.Range("K136").Value = 88239492
$ws.Range("L136").Value = 6333
$ws.Range("M136").Value = -88236942
$ws.Range("N136").Value = -11433

# WVR row 132
$ws = $wb.Worksheets.Item(8)
$ws.Range("H132").Value = 1338.8422
$ws.Range("I132").Value = 1199.6296
$ws.Range("J132").Value = 1680.5454
$ws.Range("K132").Value = 3598.8888
$ws.Range("L132").Value = 5041.6362
$ws.Range("M132").Value = -1068.8888
$ws.Range("N132").Value = -10101.6362

# WVR row 136
$ws = $wb.Worksheets.Item(8)
$ws.Range("H136").Value = 1392.0714
$ws.Range("I136").Value = 1220.421
$ws.Range("J136").Value = 1754.4445
$ws.Range("K136").Value = 3661.263
$ws.Range("L136").Value = 5263.333500000001
$ws.Range("M136").Value = -1111.263
$ws.Range("N136").Value = -10363.3335
